# Update TPM-derived NATMI metrics for Jag2-Notch3 (per updated pipeline output)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.93194
$ws.Range("H2").Value = 47.79582
$ws.Range("I2").Value = 0.9552847657129105
$ws.Range("J2").Value = 0.9552847657129107
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723263
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 51.63681424007333
$ws.Range("R2").Value = 464.73132816066
$ws.Range("S2").Value = 0.02375084273799995
$ws.Range("T2").Value = 0.02375084273799996

# Row 3
$ws.Range("G3").Value = 15.93194
$ws.Range("H3").Value = 47.79582
$ws.Range("I3").Value = 0.9552847657129105
$ws.Range("J3").Value = 0.9552847657129107
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 43.53416943745999
$ws.Range("R3").Value = 391.8075249371399
$ws.Range("S3").Value = 0.02002395436773731
$ws.Range("T3").Value = 0.02002395436773731

# Row 4
$ws.Range("G4").Value = 15.93194
$ws.Range("H4").Value = 47.79582
$ws.Range("I4").Value = 0.9552847657129105
$ws.Range("J4").Value = 0.9552847657129107
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 1981.717930860553
$ws.Range("R4").Value = 17835.46137774498
$ws.Range("S4").Value = 0.9115099686071731
$ws.Range("T4").Value = 0.9115099686071734

# Row 5
$ws.Range("I5").Value = 0.004609931913019111
$ws.Range("J5").Value = 0.004609931913019112
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723263
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 0.2491845430763333
$ws.Range("R5").Value = 2.242660887687
$ws.Range("S5").Value = 0.0001146147953247157
$ws.Range("T5").Value = 0.0001146147953247157

# Row 6
$ws.Range("I6").Value = 0.004609931913019111
$ws.Range("J6").Value = 0.004609931913019112
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("S6").Value = 0.00009662989464275835
$ws.Range("T6").Value = 0.00009662989464275837

# Row 7
$ws.Range("I7").Value = 0.004609931913019111
$ws.Range("J7").Value = 0.004609931913019112
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 9.563205716212334
$ws.Range("R7").Value = 86.06885144591101
$ws.Range("S7").Value = 0.004398687223051636
$ws.Range("T7").Value = 0.004398687223051638

# Row 8
$ws.Range("G8").Value = 0.6688636666666667
$ws.Range("H8").Value = 2.006591
$ws.Range("I8").Value = 0.04010530237407027
$ws.Range("J8").Value = 0.04010530237407027
$ws.Range("M8").Value = 3.241087666666667
$ws.Range("N8").Value = 9.723263
$ws.Range("O8").Value = 0.02486257877280725
$ws.Range("P8").Value = 0.02486257877280725
$ws.Range("Q8").Value = 2.167845780714778
$ws.Range("R8").Value = 19.510612026433
$ws.Range("S8").Value = 0.0009971212394825756
$ws.Range("T8").Value = 0.0009971212394825758

# Row 9
$ws.Range("G9").Value = 0.6688636666666667
$ws.Range("H9").Value = 2.006591
$ws.Range("I9").Value = 0.04010530237407027
$ws.Range("J9").Value = 0.04010530237407027
$ws.Range("O9").Value = 0.02096124117795788
$ws.Range("P9").Value = 0.02096124117795788
$ws.Range("Q9").Value = 1.827675988939667
$ws.Range("R9").Value = 16.449083900457
$ws.Range("S9").Value = 0.0008406569155778138
$ws.Range("T9").Value = 0.0008406569155778138

# Row 10
$ws.Range("G10").Value = 0.6688636666666667
$ws.Range("H10").Value = 2.006591
$ws.Range("I10").Value = 0.04010530237407027
$ws.Range("J10").Value = 0.04010530237407027
$ws.Range("M10").Value = 124.3864796666667
$ws.Range("N10").Value = 373.159439
$ws.Range("O10").Value = 0.9541761800492348
$ws.Range("P10").Value = 0.9541761800492349
$ws.Range("Q10").Value = 83.19759687360546
$ws.Range("R10").Value = 748.7783718624491
$ws.Range("S10").Value = 0.03826752421900988
$ws.Range("T10").Value = 0.03826752421900988
